$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of row label (first cell text, trimmed) -> value to place in the
# second (currently empty) cell of that row.
$values = @{
    "Ratio" = "0.5714";
    "Answer Recall Lenient (ARL)" = "0.4444";
    "Answer Recall Strict (ARS)" = "0.3333";
    "Answer Recall Average (ARA)" = "0.3888";
}

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    if ($row.Cells.Count -ne 2) {
        continue
    }

    $label = $row.Cells.Item(1).Range.Text
    $label = $label -replace "[`r`a`n]", ""
    $label = $label.Trim()

    if ($values.ContainsKey($label)) {
        $cell = $row.Cells.Item(2)
        $cell.Range.Text = $values[$label]
        $cell.Range.Font.Bold = $true
        $cell.Range.Font.Size = 12
        $cell.Range.Font.SizeBi = 12
    }
}
